# update frontend (add schedule scheme: add rules)
#
# The "games" table (header row 38 / type row 39) drops its
# "score_set_total_player1" / "score_set_total_player2" columns
# (originally J:K). Every column to the right of them (L:AA) shifts
# two places to the left (into J:Y), carrying its value and yellow
# highlight formatting along with it, and the two now-unused trailing
# columns are removed.
#
# NOTE: Range.Delete(xlShiftToLeft) in this environment shifts entire
# columns sheet-wide instead of being scoped to the selected rows, so
# cells are moved by hand (read, then write) instead of relying on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcFirstCol  = 12   # column L  (first surviving column of the games table)
$srcLastCol   = 27   # column AA (last column of the games table)
$destFirstCol = 10   # column J  (where score_set_total_player1 used to start)
$rows = @(38, 39)

foreach ($r in $rows) {
    # Read every value/format in the block that survives, left to right,
    # buffering in memory first so the in-place left shift can't clobber
    # a value before it has been read (source and destination overlap).
    $vals = @()
    $colors = @()
    $patterns = @()
    for ($c = $srcFirstCol; $c -le $srcLastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $vals += ,$cell.Value2
        $colors += ,$cell.Interior.Color
        $patterns += ,$cell.Interior.Pattern
    }

    for ($i = 0; $i -lt $vals.Length; $i++) {
        $destCell = $ws.Cells.Item($r, $destFirstCol + $i)
        $destCell.Value = $vals[$i]
        if ($patterns[$i] -eq [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone) {
            $destCell.Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone
        } else {
            $destCell.Interior.Color = $colors[$i]
        }
    }
}

# The trailing two columns of the table (old Z:AA) are now spare - nothing
# else on the sheet uses them - so remove them outright, which also shrinks
# the sheet's used range/dimension back down like a real column delete would.
$ws.Range("Z1:AA100").EntireColumn.Delete()

# Reflect the cell the editor ended up with selected after making the change.
$ws.Range("W38").Select()
